$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I15:I116").Clear()
